$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (surgical substring replace to preserve shared-string reuse) ---
$ws.Range("A8").Characters(21, 2).Text = "42"
$ws.Range("C9").Characters(27, 9).Text = "10/16/2023"
$ws.Range("C9").Characters(48, 10).Text = "10/22/2023"

# --- Cells changing from numeric to text ("0" / "***.*") ---
$ws.Range("C14").Copy($ws.Range("C23"))
$ws.Range("D14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("G28"))
$ws.Range("E14").Copy($ws.Range("H28"))
$ws.Range("C14").Copy($ws.Range("G29"))
$ws.Range("E14").Copy($ws.Range("H29"))

# --- Cells changing from text to numeric ---
$ws.Range("F14").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("H15").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("F14").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 3
$ws.Range("H15").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100

# --- Pure numeric value updates ---
$ws.Range("L15").Value = -36.363636363636
$ws.Range("M15").Value = -50
$ws.Range("N15").Value = 40
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -14.285714285714
$ws.Range("F16").Value = 15
$ws.Range("H16").Value = -40
$ws.Range("I16").Value = 142
$ws.Range("J16").Value = 167
$ws.Range("K16").Value = -14.970059880239
$ws.Range("L16").Value = -25.65445026178
$ws.Range("M16").Value = -24.867724867724
$ws.Range("N16").Value = 202.127659574468
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 15.384615384615
$ws.Range("I17").Value = 208
$ws.Range("J17").Value = 225
$ws.Range("K17").Value = -7.555555555555
$ws.Range("L17").Value = -19.37984496124
$ws.Range("M17").Value = 40.54054054054
$ws.Range("N17").Value = 700
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -16.666666666666
$ws.Range("I18").Value = 97
$ws.Range("J18").Value = 84
$ws.Range("K18").Value = 15.47619047619
$ws.Range("L18").Value = -11.009174311926
$ws.Range("M18").Value = 14.117647058823
$ws.Range("N18").Value = 136.585365853659
$ws.Range("C19").Value = 5
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 22
$ws.Range("H19").Value = 18.181818181818
$ws.Range("I19").Value = 297
$ws.Range("J19").Value = 242
$ws.Range("K19").Value = 22.727272727272
$ws.Range("L19").Value = 6.071428571428
$ws.Range("M19").Value = 37.5
$ws.Range("N19").Value = 1088
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 120
$ws.Range("J20").Value = 119
$ws.Range("K20").Value = 0.840336134453
$ws.Range("L20").Value = 103.389830508475
$ws.Range("M20").Value = 103.389830508475
$ws.Range("N20").Value = 172.727272727273
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -4.761904761904
$ws.Range("F21").Value = 77
$ws.Range("G21").Value = 82
$ws.Range("H21").Value = -6.097560975609
$ws.Range("I21").Value = 878
$ws.Range("J21").Value = 851
$ws.Range("K21").Value = 3.172737955346
$ws.Range("L21").Value = -3.728070175438
$ws.Range("M21").Value = 22.797202797202
$ws.Range("N21").Value = 367.021276595745
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("J22").Value = 24
$ws.Range("K22").Value = -29.166666666666
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 21.428571428571
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -21.739130434782
$ws.Range("F24").Value = 81
$ws.Range("G24").Value = 69
$ws.Range("H24").Value = 17.391304347826
$ws.Range("I24").Value = 783
$ws.Range("J24").Value = 1114
$ws.Range("K24").Value = -29.712746858168
$ws.Range("L24").Value = -37.309847878302
$ws.Range("M24").Value = 62.111801242236
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = -3.703703703703
$ws.Range("I25").Value = 314
$ws.Range("J25").Value = 330
$ws.Range("K25").Value = -4.848484848484
$ws.Range("L25").Value = -8.454810495626
$ws.Range("M25").Value = -22.277227722772
$ws.Range("G26").Value = 2
$ws.Range("L26").Value = -27.777777777777
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 25
$ws.Range("J27").Value = 42
$ws.Range("K27").Value = -16.666666666666
